$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(9, 8).Value = 70
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 13).ClearContents()
$ws.Cells.Item(64, 8).Value = 4352.9414
$ws.Cells.Item(64, 10).Value = 4000
$ws.Cells.Item(64, 12).Value = 4000
$ws.Cells.Item(64, 14).Value = -4496
$ws.Cells.Item(67, 8).Value = 4352.9414
$ws.Cells.Item(67, 10).Value = 4000
$ws.Cells.Item(67, 12).Value = 4000
$ws.Cells.Item(67, 14).Value = -5716
$ws.Cells.Item(112, 8).Value = 1697
$ws.Cells.Item(112, 10).Value = 1883.0667
$ws.Cells.Item(112, 12).Value = 5649.2001
$ws.Cells.Item(112, 14).Value = -7865.2001
$ws.Cells.Item(138, 8).Value = 2914.45
$ws.Cells.Item(138, 10).Value = 3081.4285
$ws.Cells.Item(138, 12).Value = 9244.2855
$ws.Cells.Item(138, 14).Value = -19524.2855
# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 1479.3529
$ws.Cells.Item(2, 9).Value = 1603.5714
$ws.Cells.Item(2, 10).Value = 899.6667
$ws.Cells.Item(2, 11).Value = 1603.5714
$ws.Cells.Item(2, 12).Value = 899.6667
$ws.Cells.Item(2, 13).Value = -1490.5714
$ws.Cells.Item(2, 14).Value = -1125.6667
$ws.Cells.Item(32, 8).Value = 4629.3657
$ws.Cells.Item(32, 9).Value = 3092.7715
$ws.Cells.Item(32, 10).Value = 13592.833
$ws.Cells.Item(32, 11).Value = 3092.7715
$ws.Cells.Item(32, 12).Value = 13592.833
$ws.Cells.Item(32, 13).Value = -2805.7715
$ws.Cells.Item(32, 14).Value = -14166.833
$ws.Cells.Item(102, 8).Value = 1577.75
$ws.Cells.Item(102, 9).Value = 1577.75
$ws.Cells.Item(102, 11).Value = 1577.75
$ws.Cells.Item(102, 13).Value = 44.25
$ws.Cells.Item(116, 8).Value = 1479.3529
$ws.Cells.Item(116, 9).Value = 1603.5714
$ws.Cells.Item(116, 10).Value = 899.6667
$ws.Cells.Item(116, 11).Value = 1603.5714
$ws.Cells.Item(116, 12).Value = 899.6667
$ws.Cells.Item(116, 13).Value = 690.4286
$ws.Cells.Item(116, 14).Value = -5487.6667
$ws.Cells.Item(132, 8).Value = 3337.9546
$ws.Cells.Item(132, 9).Value = 2717.0715
$ws.Cells.Item(132, 11).Value = 8151.2145
$ws.Cells.Item(132, 13).Value = -5621.2145
# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 1479.3529
$ws.Cells.Item(3, 9).Value = 1603.5714
$ws.Cells.Item(3, 10).Value = 899.6667
$ws.Cells.Item(3, 11).Value = 1603.5714
$ws.Cells.Item(3, 12).Value = 899.6667
$ws.Cells.Item(3, 13).Value = -1489.5714
$ws.Cells.Item(3, 14).Value = -1127.6667
$ws.Cells.Item(81, 8).Value = 74926.664
$ws.Cells.Item(81, 10).Value = 74926.664
$ws.Cells.Item(81, 12).Value = 74926.664
$ws.Cells.Item(81, 14).Value = -77048.664
$ws.Cells.Item(84, 8).Value = 74926.664
$ws.Cells.Item(84, 10).Value = 74926.664
$ws.Cells.Item(84, 12).Value = 224779.992
$ws.Cells.Item(84, 14).Value = -235387.992
$ws.Cells.Item(94, 8).Value = 1099.6154
$ws.Cells.Item(94, 9).Value = 449.8889
$ws.Cells.Item(94, 11).Value = 449.8889
$ws.Cells.Item(94, 13).Value = 1.111100000000022
$ws.Cells.Item(105, 8).Value = 3328.4443
$ws.Cells.Item(105, 9).Value = 2807
$ws.Cells.Item(105, 10).Value = 7500
$ws.Cells.Item(105, 11).Value = 2807
$ws.Cells.Item(105, 12).Value = 7500
$ws.Cells.Item(105, 13).Value = -1060
$ws.Cells.Item(105, 14).Value = -10994
$ws.Cells.Item(107, 8).Value = 250
$ws.Cells.Item(107, 9).Value = 250
$ws.Cells.Item(107, 11).Value = 250
$ws.Cells.Item(107, 13).Value = 1670
$ws.Cells.Item(134, 8).Value = 2125.5715
$ws.Cells.Item(134, 9).Value = 1816.2
$ws.Cells.Item(134, 10).Value = 2899
$ws.Cells.Item(134, 11).Value = 5448.6
$ws.Cells.Item(134, 12).Value = 8697
$ws.Cells.Item(134, 13).Value = -2913.6
$ws.Cells.Item(134, 14).Value = -13767
# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = 2227.3333
$ws.Cells.Item(16, 9).Value = 2227.3333
$ws.Cells.Item(16, 11).Value = 2227.3333
$ws.Cells.Item(16, 13).Value = -1940.3333
$ws.Cells.Item(60, 8).Value = 34426.406
$ws.Cells.Item(60, 9).Value = 15333.333
$ws.Cells.Item(60, 10).Value = 36813.043
$ws.Cells.Item(60, 11).Value = 15333.333
$ws.Cells.Item(60, 12).Value = 36813.043
$ws.Cells.Item(60, 13).Value = -14822.333
$ws.Cells.Item(60, 14).Value = -37835.043
$ws.Cells.Item(62, 8).Value = 4933.3335
$ws.Cells.Item(62, 10).Value = 4900
$ws.Cells.Item(62, 12).Value = 4900
$ws.Cells.Item(62, 14).Value = -6148
$ws.Cells.Item(65, 8).Value = 4933.3335
$ws.Cells.Item(65, 10).Value = 4900
$ws.Cells.Item(65, 12).Value = 24500
$ws.Cells.Item(65, 14).Value = -30740
$ws.Cells.Item(86, 8).Value = 8665.571
$ws.Cells.Item(86, 9).Value = 10855.286
$ws.Cells.Item(86, 10).Value = 6475.857
$ws.Cells.Item(86, 11).Value = 10855.286
$ws.Cells.Item(86, 12).Value = 6475.857
$ws.Cells.Item(86, 13).Value = -9732.286
$ws.Cells.Item(86, 14).Value = -8721.857
$ws.Cells.Item(89, 8).Value = 8665.571
$ws.Cells.Item(89, 9).Value = 10855.286
$ws.Cells.Item(89, 10).Value = 6475.857
$ws.Cells.Item(89, 11).Value = 54276.43
$ws.Cells.Item(89, 12).Value = 32379.285
$ws.Cells.Item(89, 13).Value = -48660.43
$ws.Cells.Item(89, 14).Value = -43611.285
$ws.Cells.Item(94, 8).Value = 1137.25
$ws.Cells.Item(94, 9).Value = 1125
$ws.Cells.Item(94, 10).Value = 1149.5
$ws.Cells.Item(94, 11).Value = 1125
$ws.Cells.Item(94, 12).Value = 1149.5
$ws.Cells.Item(94, 13).Value = -674
$ws.Cells.Item(94, 14).Value = -2051.5
$ws.Cells.Item(113, 8).Value = 2227.3333
$ws.Cells.Item(113, 9).Value = 2227.3333
$ws.Cells.Item(113, 11).Value = 2227.3333
$ws.Cells.Item(113, 13).Value = -57.33329999999978
$ws.Cells.Item(134, 8).Value = 5000
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 13).ClearContents()
$ws.Cells.Item(141, 8).Value = 69885.75
$ws.Cells.Item(141, 10).Value = 59847.668
$ws.Cells.Item(141, 12).Value = 59847.668
$ws.Cells.Item(141, 14).Value = -70207.66800000001
# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(18, 8).Value = 665
$ws.Cells.Item(18, 9).Value = 665
$ws.Cells.Item(18, 11).Value = 1995
$ws.Cells.Item(18, 13).Value = -1826
$ws.Cells.Item(40, 8).Value = 576.25
$ws.Cells.Item(40, 9).Value = 576.25
$ws.Cells.Item(40, 11).Value = 2305
$ws.Cells.Item(40, 13).Value = -2236
$ws.Cells.Item(56, 8).Value = 18478.088
$ws.Cells.Item(56, 9).Value = 18478.088
$ws.Cells.Item(56, 11).Value = 18478.088
$ws.Cells.Item(56, 13).Value = -17948.088
# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(113, 8).Value = 1949
$ws.Cells.Item(113, 9).Value = 1949
$ws.Cells.Item(113, 11).Value = 1949
$ws.Cells.Item(113, 13).Value = 221
$ws.Cells.Item(122, 8).Value = 997
$ws.Cells.Item(122, 9).Value = 997
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 2991
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(122, 14).Value = -541
$ws.Cells.Item(139, 8).Value = 74193.60000000001
$ws.Cells.Item(139, 10).Value = 74193.60000000001
$ws.Cells.Item(139, 12).Value = 74193.60000000001
$ws.Cells.Item(139, 14).Value = -84473.60000000001
# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(42, 8).Value = 22306
$ws.Cells.Item(42, 10).Value = 22306
$ws.Cells.Item(42, 12).Value = 22306
$ws.Cells.Item(42, 14).Value = -23432
$ws.Cells.Item(49, 8).Value = 22306
$ws.Cells.Item(49, 10).Value = 22306
$ws.Cells.Item(49, 12).Value = 22306
$ws.Cells.Item(49, 14).Value = -22600
$ws.Cells.Item(122, 8).Value = 10000
$ws.Cells.Item(122, 9).Value = 10000
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 30000
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(122, 14).Value = -27550
# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value = 14165.5
$ws.Cells.Item(81, 9).Value = 11749.75
$ws.Cells.Item(81, 11).Value = 23499.5
$ws.Cells.Item(81, 13).Value = -22438.5
$ws.Cells.Item(84, 8).Value = 14165.5
$ws.Cells.Item(84, 9).Value = 11749.75
$ws.Cells.Item(84, 11).Value = 117497.5
$ws.Cells.Item(84, 13).Value = -112193.5
$ws.Cells.Item(96, 8).Value = 3000
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 799.5
$ws.Cells.Item(107, 10).Value = 799
$ws.Cells.Item(107, 12).Value = 2397
$ws.Cells.Item(107, 14).Value = -6237
$ws.Cells.Item(132, 8).Value = 3909.818
$ws.Cells.Item(132, 9).Value = 3445.7273
$ws.Cells.Item(132, 11).Value = 10337.1819
$ws.Cells.Item(132, 13).Value = -7807.1819
$ws.Cells.Item(136, 8).Value = 2262.6667
$ws.Cells.Item(136, 9).Value = 1711.6875
$ws.Cells.Item(136, 10).Value = 3364.625
$ws.Cells.Item(136, 11).Value = 5135.0625
$ws.Cells.Item(136, 12).Value = 10093.875
$ws.Cells.Item(136, 13).Value = -2585.0625
$ws.Cells.Item(136, 14).Value = -15193.875
